# Auto-generated edit script: update cryptos list prices/volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.482.96'
$ws.Range("E2").Value = '  +0.73%  '
$ws.Range("D3").Value = '2.426.66'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'565.63"
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("E10").Value = '  +0.37%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = '  +1.54%  '
$ws.Range("D13").Value = "'26.82"
$ws.Range("E13").Value = '  +4.85%  '
$ws.Range("E14").Value = '  +4.01%  '
$ws.Range("D15").Value = '2.865.31'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = '62.310.38'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '2.423.59'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = "'6.96"
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("D20").Value = "'323.79"
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").Value = "'67.00"
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("E24").Value = '  +4.98%  '
$ws.Range("D25").Value = "'596.26"
$ws.Range("E25").Value = '  +5.81%  '
$ws.Range("D26").Value = "'8.58"
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("E27").Value = '  +6.89%  '
$ws.Range("D28").Value = '2.544.57'
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = "'8.42"
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("E33").Value = '  -0.51%  '
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("D35").Value = "'4.85"
$ws.Range("E35").Value = '  +1.49%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("D39").Value = "'5.34"
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("D40").Value = "'147.40"
$ws.Range("E40").Value = '  -3.34%  '
$ws.Range("E41").Value = '  +1.32%  '
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  +9.29%  '
$ws.Range("D44").Value = "'148.04"
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("D45").Value = "'3.67"
$ws.Range("D46").Value = "'0.0535"
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").Value = "'20.49"
$ws.Range("E47").Value = '  +3.09%  '
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").Value = "'1.10"
$ws.Range("E51").Value = '  +4.29%  '
